$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 379.27274
$ws.Range("I5").Value = 582
$ws.Range("J5").Value = 136
$ws.Range("K5").Value = 582
$ws.Range("L5").Value = 136
$ws.Range("M5").Value = -467
$ws.Range("N5").Value = -366

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 107.63636
$ws.Range("I38").Value = 20.555555
$ws.Range("J38").Value = 499.5
$ws.Range("K38").Value = 61.66666499999999
$ws.Range("L38").Value = 1498.5
$ws.Range("M38").Value = 310.333335
$ws.Range("N38").Value = -2242.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3750
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 3750
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 11250
$ws.Range("M70").Value = $null
$ws.Range("N70").Value = -11790

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 3750
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 3750
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 11250
$ws.Range("M73").Value = $null
$ws.Range("N73").Value = -13122

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I80").Value = 775
$ws.Range("J80").Value = 2003
$ws.Range("K80").Value = 2325
$ws.Range("L80").Value = 6009
$ws.Range("M80").Value = -1327
$ws.Range("N80").Value = -8005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I83").Value = 775
$ws.Range("J83").Value = 2003
$ws.Range("K83").Value = 6975
$ws.Range("L83").Value = 18027
$ws.Range("M83").Value = -1983
$ws.Range("N83").Value = -28011

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3760.6
$ws.Range("I100").Value = 3400
$ws.Range("K100").Value = 3400
$ws.Range("M100").Value = -2859

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4000
$ws.Range("I61").Value = 4000
$ws.Range("K61").Value = 4000
$ws.Range("M61").Value = -3788

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1914
$ws.Range("I74").Value = 1871
$ws.Range("K74").Value = 1871
$ws.Range("M74").Value = -997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1914
$ws.Range("I77").Value = 1871
$ws.Range("K77").Value = 9355
$ws.Range("M77").Value = -4987

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 10552446
$ws.Range("I102").Value = 787779.9
$ws.Range("K102").Value = 787779.9
$ws.Range("M102").Value = -786157.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4000
$ws.Range("I136").Value = 4000
$ws.Range("K136").Value = 12000
$ws.Range("M136").Value = -9450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 111498.5
$ws.Range("I94").Value = 138879.38
$ws.Range("J94").Value = 1975
$ws.Range("K94").Value = 138879.38
$ws.Range("L94").Value = 1975
$ws.Range("M94").Value = -138428.38
$ws.Range("N94").Value = -2877

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 205
$ws.Range("I7").Value = 234.17647
$ws.Range("K7").Value = 234.17647
$ws.Range("M7").Value = -121.17647

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 4000000
$ws.Range("I23").Value = 4000000
$ws.Range("K23").Value = 4000000
$ws.Range("M23").Value = -3999760

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H27").Value = 4000000
$ws.Range("I27").Value = 4000000
$ws.Range("K27").Value = 4000000
$ws.Range("M27").Value = -3999808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 12892.5
$ws.Range("I41").Value = 2190
$ws.Range("K41").Value = 2190
$ws.Range("M41").Value = -1762

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4000
$ws.Range("I58").Value = 4000
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 4000
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -3797
$ws.Range("N58").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1102.5
$ws.Range("I105").Value = 803.3333
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 803.3333
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = 943.6667
$ws.Range("N105").Value = -5494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 857.7692
$ws.Range("I107").Value = 795.8889
$ws.Range("J107").Value = 997
$ws.Range("K107").Value = 795.8889
$ws.Range("L107").Value = 997
$ws.Range("M107").Value = 1124.1111
$ws.Range("N107").Value = -4837

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3766
$ws.Range("I122").Value = 957.625
$ws.Range("J122").Value = 14999.5
$ws.Range("K122").Value = 2872.875
$ws.Range("L122").Value = 44998.5
$ws.Range("M122").Value = -422.875
$ws.Range("N122").Value = -49898.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4000
$ws.Range("I136").Value = 4000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 12000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -9450
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 3401
$ws.Range("I112").Value = 1201.3334
$ws.Range("J112").Value = 10000
$ws.Range("K112").Value = 3604.0002
$ws.Range("L112").Value = 30000
$ws.Range("M112").Value = -2496.0002
$ws.Range("N112").Value = -32216

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2588.625
$ws.Range("I139").Value = 2387
$ws.Range("J139").Value = 4000
$ws.Range("K139").Value = 7161
$ws.Range("L139").Value = 12000
$ws.Range("M139").Value = -2021
$ws.Range("N139").Value = -22280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 4528
$ws.Range("I22").Value = 1585
$ws.Range("J22").Value = 5999.5
$ws.Range("K22").Value = 1585
$ws.Range("L22").Value = 5999.5
$ws.Range("M22").Value = -1056
$ws.Range("N22").Value = -7057.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5879.923
$ws.Range("I122").Value = 4493.3335
$ws.Range("K122").Value = 13480.0005
$ws.Range("M122").Value = -11030.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2068.1936
$ws.Range("I22").Value = 1902.0526
$ws.Range("J22").Value = 2331.25
$ws.Range("K22").Value = 1902.0526
$ws.Range("L22").Value = 2331.25
$ws.Range("M22").Value = -1607.0526
$ws.Range("N22").Value = -2921.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2068.1936
$ws.Range("I27").Value = 1902.0526
$ws.Range("J27").Value = 2331.25
$ws.Range("K27").Value = 1902.0526
$ws.Range("L27").Value = 2331.25
$ws.Range("M27").Value = -1795.0526
$ws.Range("N27").Value = -2545.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1092.8889
$ws.Range("I55").Value = 660.1667
$ws.Range("J55").Value = 1958.3334
$ws.Range("K55").Value = 660.1667
$ws.Range("L55").Value = 1958.3334
$ws.Range("M55").Value = -487.1667
$ws.Range("N55").Value = -2304.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 333333340
$ws.Range("I93").Value = 333333340
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 333333340
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -333332092
$ws.Range("N93").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3585.8928
$ws.Range("I122").Value = 3023.4614
$ws.Range("K122").Value = 9070.3842
$ws.Range("M122").Value = -6620.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2183
$ws.Range("I122").Value = 2150
$ws.Range("J122").Value = 2199.5
$ws.Range("K122").Value = 2150
$ws.Range("L122").Value = 6598.5
$ws.Range("M122").Value = -4000
$ws.Range("N122").Value = -11498.5
